$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows before row 169; this shifts the existing rows 169-186
# down to 171-188, preserving all of their data and formatting.
$ws.Rows("169:170").Insert()

# Populate the two newly inserted rows with the new weekly price entries.
# Row 169 - "$/caja 18 kilos" entry
$ws.Range("A169").Value = 9
$ws.Range("B169").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C169").Value = "Metropolitana"
$ws.Range("D169").Value = 45194
$ws.Range("D169").NumberFormat = $ws.Range("D171").NumberFormat()
$ws.Range("E169").Value = 13
$ws.Range("F169").Value = 100114002
$ws.Range("G169").Value = "Camote"
$ws.Range("H169").Value = "Sin especificar"
$ws.Range("I169").Value = "Primera"
$ws.Range("J169").Value = 700
$ws.Range("K169").Value = 17000
$ws.Range("L169").Value = 18000
$ws.Range("M169").Value = 17500
$ws.Range("N169").Value = "$/caja 18 kilos"
$ws.Range("O169").Value = "Perú"
$ws.Range("P169").Value = 972
$ws.Range("Q169").Value = 18
$ws.Range("R169").Value = "Hortaliza"

# Row 170 - "$/malla 18 kilos" entry
$ws.Range("A170").Value = 9
$ws.Range("B170").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C170").Value = "Metropolitana"
$ws.Range("D170").Value = 45194
$ws.Range("D170").NumberFormat = $ws.Range("D171").NumberFormat()
$ws.Range("E170").Value = 13
$ws.Range("F170").Value = 100114002
$ws.Range("G170").Value = "Camote"
$ws.Range("H170").Value = "Sin especificar"
$ws.Range("I170").Value = "Primera"
$ws.Range("J170").Value = 880
$ws.Range("K170").Value = 15000
$ws.Range("L170").Value = 16000
$ws.Range("M170").Value = 15500
$ws.Range("N170").Value = "$/malla 18 kilos"
$ws.Range("O170").Value = "Perú"
$ws.Range("P170").Value = 861
$ws.Range("Q170").Value = 18
$ws.Range("R170").Value = "Hortaliza"
